$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# `cryptos` price/volume cells are stored as plain text (e.g. "212.24",
# "26.237.91", "  -0.31%  ") even when the text happens to look like a
# number. Writing such a string straight into `.Value` lets Excel's normal
# auto-detection kick in and silently reinterpret it as a real number
# (dropping trailing zeros / introducing float noise, e.g. "63.80" -> 63.8,
# "1.00" -> 1). To avoid that, cells whose new text is numeric-looking are
# switched to Text format before the write, then restored to the default
# "Normal" style afterwards (the value itself stays text; only the leftover
# number-format flag is cleaned up so the cell style matches the original).
function Set-CellText($ref, $text) {
    $cell = $ws.Range($ref)
    $looksNumeric = $text -match '^[+-]?(\d+\.?\d*|\.\d+)$'
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-CellText "D2" "26.237.91"
Set-CellText "E2" "  -0.31%  "
Set-CellText "D3" "1.591.21"
Set-CellText "E3" "  +0.04%  "
Set-CellText "E4" "  +0.05%  "
Set-CellText "D5" "212.24"
Set-CellText "E5" "  +0.44%  "
Set-CellText "E6" "  -0.71%  "
Set-CellText "E7" "  +0.09%  "
Set-CellText "E8" "  -0.75%  "
Set-CellText "E9" "  -0.58%  "
Set-CellText "D10" "19.01"
Set-CellText "E10" "  -2.03%  "
Set-CellText "D13" "1.592.85"
Set-CellText "E13" "  +0.38%  "
Set-CellText "E14" "  -1.79%  "
Set-CellText "D15" "0.509"
Set-CellText "E15" "  -2.54%  "
Set-CellText "D16" "63.80"
Set-CellText "E16" "  -1.20%  "
Set-CellText "D17" "26.234.14"
Set-CellText "E17" "  -0.34%  "
Set-CellText "E18" "  -0.74%  "
Set-CellText "D19" "215.43"
Set-CellText "E19" "  +1.63%  "
Set-CellText "D20" "7.31"
Set-CellText "E20" "  -2.80%  "
Set-CellText "E21" "  +0.02%  "
Set-CellText "D22" "4.29"
Set-CellText "E22" "  -0.11%  "
Set-CellText "D23" "9.04"
Set-CellText "E23" "  +0.35%  "
Set-CellText "E24" "  -1.37%  "
Set-CellText "D25" "144.71"
Set-CellText "E25" "  +0.42%  "
Set-CellText "E26" "  +0.03%  "
Set-CellText "D27" "6.96"
Set-CellText "E27" "  -1.45%  "
Set-CellText "E28" "  -0.89%  "
Set-CellText "D29" "15.11"
Set-CellText "E29" "  -0.92%  "
Set-CellText "E30" "  -2.80%  "
Set-CellText "E31" "  +0.27%  "
Set-CellText "D32" "3.18"
Set-CellText "E32" "  -1.13%  "
Set-CellText "D33" "1.412.14"
Set-CellText "E33" "  +6.22%  "
Set-CellText "D34" "2.95"
Set-CellText "E34" "  -1.38%  "
Set-CellText "E35" "  -0.40%  "
Set-CellText "D36" "1.46"
Set-CellText "E36" "  -0.99%  "
Set-CellText "D37" "0.580"
Set-CellText "E37" "  -3.78%  "
Set-CellText "D39" "0.822"
Set-CellText "E39" "  +0.43%  "
Set-CellText "D40" "5.84"
Set-CellText "E40" "  +2.24%  "
Set-CellText "E41" "  -0.02%  "
Set-CellText "D42" "0.978"
Set-CellText "E42" "  -1.12%  "
Set-CellText "E43" "  +0.09%  "
Set-CellText "E44" "  +0.05%  "
Set-CellText "D45" "1.727.89"
Set-CellText "E45" "  +0.13%  "
Set-CellText "D46" "60.85"
Set-CellText "E46" "  -1.66%  "
Set-CellText "D47" "86.12"
Set-CellText "E47" "  -2.26%  "
Set-CellText "E48" "  -0.01%  "
Set-CellText "E49" "  -0.66%  "
Set-CellText "E50" "  -2.51%  "
Set-CellText "D51" "1.00"
Set-CellText "E51" "  -0.05%  "
